$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the old 3-column layout (A1:C4) entirely
$ws.Range("A1:C4").Clear()

# Write the new single-column data (A1:A8)
$values = @(
    "信息与通信工程学院",
    "通信工程",
    "电子信息工程",
    "广播电视工程（智能视听技术方向）",
    "数字媒体技术",
    "物联网工程",
    "人工智能",
    "智能装备系统（演艺工程与智能技术方向）"
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}

$win = $excel.ActiveWindow
try { $win.ScrollRow = 3 } catch {}
try { $win.ScrollColumn = 1 } catch {}

$ws.Range("E8").Select()
